$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E2" = 3
    "F2" = 1
    "G2" = 0.9919543333333333
    "H2" = 2.975863
    "I2" = 0.008811579445878926
    "J2" = 0.008811579445878926
    "M2" = 0.3259846666666666
    "N2" = 0.977954
    "O2" = 0.004039686607851542
    "P2" = 0.004039686607851542
    "Q2" = 0.3233619027002222
    "R2" = 2.910257124302
    "S2" = 0.00003559601948153701
    "T2" = 0.00003559601948153701
    "E3" = 3
    "F3" = 1
    "G3" = 0.9919543333333333
    "H3" = 2.975863
    "I3" = 0.008811579445878926
    "J3" = 0.008811579445878926
    "O3" = 0.9948000963506443
    "P3" = 0.9948000963506443
    "Q3" = 79.63005133544021
    "R3" = 716.670462018962
    "S3" = 0.008765760081761711
    "T3" = 0.008765760081761711
    "E4" = 3
    "F4" = 1
    "G4" = 0.9919543333333333
    "H4" = 2.975863
    "I4" = 0.008811579445878926
    "J4" = 0.008811579445878926
    "K4" = 1
    "L4" = 0.3333333333333333
    "M4" = 0.09362433333333332
    "N4" = 0.280873
    "O4" = 0.001160217041504085
    "P4" = 0.001160217041504085
    "Q4" = 0.09287106315544442
    "R4" = 0.835839568399
    "S4" = 0.00001022334463567585
    "T4" = 0.00001022334463567585
    "I5" = 0.6711393126876655
    "J5" = 0.6711393126876655
    "M5" = 0.3259846666666666
    "N5" = 0.977954
    "O5" = 0.004039686607851542
    "P5" = 0.004039686607851542
    "Q5" = 24.62905617097977
    "R5" = 221.661505538818
    "S5" = 0.002711192493467051
    "T5" = 0.002711192493467051
    "I6" = 0.6711393126876655
    "J6" = 0.6711393126876655
    "O6" = 0.9948000963506443
    "P6" = 0.9948000963506443
    "S6" = 0.6676494529263948
    "T6" = 0.6676494529263948
    "I7" = 0.6711393126876655
    "J7" = 0.6711393126876655
    "K7" = 1
    "L7" = 0.3333333333333333
    "M7" = 0.09362433333333332
    "N7" = 0.280873
    "O7" = 0.001160217041504085
    "P7" = 0.001160217041504085
    "Q7" = 7.073581062004554
    "R7" = 63.66222955804099
    "S7" = 0.0007786672678035685
    "T7" = 0.0007786672678035685
    "G8" = 36.02919333333333
    "H8" = 108.08758
    "I8" = 0.3200491078664556
    "J8" = 0.3200491078664556
    "M8" = 0.3259846666666666
    "N8" = 0.977954
    "O8" = 0.004039686607851542
    "P8" = 0.004039686607851542
    "Q8" = 11.74496457903556
    "R8" = 105.70468121132
    "S8" = 0.001292898094902954
    "T8" = 0.001292898094902954
    "G9" = 36.02919333333333
    "H9" = 108.08758
    "I9" = 0.3200491078664556
    "J9" = 0.3200491078664556
    "O9" = 0.9948000963506443
    "P9" = 0.9948000963506443
    "Q9" = 2892.276809827435
    "R9" = 26030.49128844692
    "S9" = 0.3183848833424878
    "T9" = 0.3183848833424878
    "G10" = 36.02919333333333
    "H10" = 108.08758
    "I10" = 0.3200491078664556
    "J10" = 0.3200491078664556
    "K10" = 1
    "L10" = 0.3333333333333333
    "M10" = 0.09362433333333332
    "N10" = 0.280873
    "O10" = 0.001160217041504085
    "P10" = 0.001160217041504085
    "Q10" = 3.373209206371111
    "R10" = 30.35888285734
    "S10" = 0.000371326429064841
    "T10" = 0.000371326429064841
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
